$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "Move to location (10, 2) and remove the toolkit."
$ws.Range("B1").Value = "['Robot42']"
$ws.Range("E1").Value = "(10, 2)"

# Row 2
$ws.Range("A2").Value = "Move to location (6, 6) and remove the liquid spill."
$ws.Range("E2").Value = "(6, 6)"

# Row 3
$ws.Range("A3").Value = "Move to location (5, 3) and remove the large debris."
$ws.Range("B3").Value = "['Robot2', 'Robot39']"
$ws.Range("E3").Value = "(5, 3)"

# Row 4
$ws.Range("A4").Value = "Move to location (6, 10) and remove the dust."
$ws.Range("E4").Value = "(6, 10)"

# Row 5
$ws.Range("A5").Value = "Move to location (4, 8) and remove the grass."
$ws.Range("B5").Value = "['Robot21']"
$ws.Range("E5").Value = "(4, 8)"

# Row 6
$ws.Range("A6").Value = "Move to location (8, 7) and remove the small debris."
$ws.Range("B6").Value = "['Robot8', 'Robot50']"
$ws.Range("E6").Value = "(8, 7)"

# Row 7
$ws.Range("A7").Value = "Move to location (1, 10) and remove the vehicle."
$ws.Range("E7").Value = "(1, 10)"

# Row 8
$ws.Range("A8").Value = "Move to location (2, 12) and remove the construction materials."
$ws.Range("B8").Value = "['Robot22', 'Robot19', 'Robot13']"
$ws.Range("E8").Value = "(2, 12)"

# Row 9
$ws.Range("A9").Value = "Move to location (8, 9) and remove the tree branches."
$ws.Range("E9").Value = "(8, 9)"

# Row 10
$ws.Range("A10").Value = "Move to location (7, 5) and remove the screws."
$ws.Range("E10").Value = "(7, 5)"
